$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 3335333
$ws.Range("I6").Value = 3335333
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 10005999
$ws.Range("L6").Value = 0
$ws.Range("M6").Value = -10005887
$ws.Range("N6").ClearContents()
$ws.Range("H8").Value = 114.71429
$ws.Range("I8").Value = 114.71429
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 344.14287
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -205.14287
$ws.Range("N8").ClearContents()
$ws.Range("H137").Value = 366137.97
$ws.Range("I137").Value = 518650.66
$ws.Range("J137").Value = 47247.816
$ws.Range("K137").Value = 1555951.98
$ws.Range("L137").Value = 141743.448
$ws.Range("M137").Value = -1553401.98
$ws.Range("N137").Value = -146843.448

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H54").Value = 34863.332
$ws.Range("J54").Value = 34863.332
$ws.Range("L54").Value = 34863.332
$ws.Range("N54").Value = -36401.332
$ws.Range("H61").Value = 5916.4546
$ws.Range("I61").Value = 7478.1333
$ws.Range("J61").Value = 2570
$ws.Range("K61").Value = 7478.1333
$ws.Range("L61").Value = 2570
$ws.Range("M61").Value = -7266.1333
$ws.Range("N61").Value = -2994
$ws.Range("H74").Value = 3537
$ws.Range("I74").Value = 795.0952
$ws.Range("J74").Value = 6567.5264
$ws.Range("K74").Value = 795.0952
$ws.Range("L74").Value = 6567.5264
$ws.Range("M74").Value = 78.90480000000002
$ws.Range("N74").Value = -8315.526399999999
$ws.Range("H77").Value = 3537
$ws.Range("I77").Value = 795.0952
$ws.Range("J77").Value = 6567.5264
$ws.Range("K77").Value = 3975.476
$ws.Range("L77").Value = 32837.632
$ws.Range("M77").Value = 392.5240000000003
$ws.Range("N77").Value = -41573.632
$ws.Range("H110").Value = 765.8570999999999
$ws.Range("I110").Value = 635.375
$ws.Range("J110").Value = 939.8333
$ws.Range("K110").Value = 635.375
$ws.Range("L110").Value = 939.8333
$ws.Range("M110").Value = 1409.625
$ws.Range("N110").Value = -5029.8333
$ws.Range("H132").Value = 3907833.5
$ws.Range("I132").Value = 5209258
$ws.Range("J132").Value = 3559.625
$ws.Range("K132").Value = 15627774
$ws.Range("L132").Value = 10678.875
$ws.Range("M132").Value = -15625244
$ws.Range("N132").Value = -15738.875
$ws.Range("H136").Value = 5916.4546
$ws.Range("I136").Value = 7478.1333
$ws.Range("J136").Value = 2570
$ws.Range("K136").Value = 22434.3999
$ws.Range("L136").Value = 7710
$ws.Range("M136").Value = -19884.3999
$ws.Range("N136").Value = -12810

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 15899078
$ws.Range("I134").Value = 22257656
$ws.Range("K134").Value = 66772968
$ws.Range("M134").Value = -66770433

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 7751.5
$ws.Range("I3").Value = 3000
$ws.Range("J3").Value = 9335.333000000001
$ws.Range("K3").Value = 3000
$ws.Range("L3").Value = 9335.333000000001
$ws.Range("M3").Value = -2887
$ws.Range("N3").Value = -9561.333000000001
$ws.Range("H31").Value = 8998.741
$ws.Range("J31").Value = 20533.092
$ws.Range("L31").Value = 20533.092
$ws.Range("N31").Value = -21123.092
$ws.Range("H34").Value = 8998.741
$ws.Range("J34").Value = 20533.092
$ws.Range("L34").Value = 20533.092
$ws.Range("N34").Value = -20937.092
$ws.Range("H58").Value = 5332722.5
$ws.Range("I58").Value = 7193650
$ws.Range("J58").Value = 15785.714
$ws.Range("K58").Value = 7193650
$ws.Range("L58").Value = 15785.714
$ws.Range("M58").Value = -7193447
$ws.Range("N58").Value = -16191.714
$ws.Range("H86").Value = 2889
$ws.Range("I86").Value = 2972.5
$ws.Range("J86").Value = 2838.9
$ws.Range("K86").Value = 2972.5
$ws.Range("L86").Value = 2838.9
$ws.Range("M86").Value = -1849.5
$ws.Range("N86").Value = -5084.9
$ws.Range("H89").Value = 2889
$ws.Range("I89").Value = 2972.5
$ws.Range("J89").Value = 2838.9
$ws.Range("K89").Value = 14862.5
$ws.Range("L89").Value = 14194.5
$ws.Range("M89").Value = -9246.5
$ws.Range("N89").Value = -25426.5
$ws.Range("H132").Value = 6945967.5
$ws.Range("I132").Value = 10417486
$ws.Range("J132").Value = 2929.875
$ws.Range("K132").Value = 31252458
$ws.Range("L132").Value = 8789.625
$ws.Range("M132").Value = -31249928
$ws.Range("N132").Value = -13849.625
$ws.Range("H134").Value = 7623055.5
$ws.Range("I134").Value = 8334315
$ws.Range("J134").Value = 5683257
$ws.Range("K134").Value = 25002945
$ws.Range("L134").Value = 17049771
$ws.Range("M134").Value = -25000410
$ws.Range("N134").Value = -17054841
$ws.Range("H136").Value = 5332722.5
$ws.Range("I136").Value = 7193650
$ws.Range("J136").Value = 15785.714
$ws.Range("K136").Value = 21580950
$ws.Range("L136").Value = 47357.142
$ws.Range("M136").Value = -21578400
$ws.Range("N136").Value = -52457.142

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 625.65
$ws.Range("I34").Value = 123.333336
$ws.Range("J34").Value = 840.9286
$ws.Range("K34").Value = 370.000008
$ws.Range("L34").Value = 2522.7858
$ws.Range("M34").Value = -286.000008
$ws.Range("N34").Value = -2690.7858
$ws.Range("H39").Value = 2300
$ws.Range("J39").Value = 2381.818
$ws.Range("L39").Value = 7145.454000000001
$ws.Range("N39").Value = -7733.454000000001
$ws.Range("H62").Value = 2200
$ws.Range("I62").Value = 1923.0769
$ws.Range("J62").Value = 2714.2856
$ws.Range("K62").Value = 5769.2307
$ws.Range("L62").Value = 8142.8568
$ws.Range("M62").Value = -5083.2307
$ws.Range("N62").Value = -9514.856800000001
$ws.Range("H65").Value = 2200
$ws.Range("I65").Value = 1923.0769
$ws.Range("J65").Value = 2714.2856
$ws.Range("K65").Value = 17307.6921
$ws.Range("L65").Value = 24428.5704
$ws.Range("M65").Value = -13875.6921
$ws.Range("N65").Value = -31292.5704
$ws.Range("H122").Value = 595.0454999999999
$ws.Range("J122").Value = 872
$ws.Range("L122").Value = 7848
$ws.Range("N122").Value = -12748
$ws.Range("H131").Value = 22572268
$ws.Range("I131").Value = 166667070
$ws.Range("J131").Value = 1987296.2
$ws.Range("K131").Value = 500001210
$ws.Range("L131").Value = 5961888.6
$ws.Range("M131").Value = -499996170
$ws.Range("N131").Value = -5971968.6

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 125000664
$ws.Range("I122").Value = 125000664
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 375001992
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -374999542
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 181826740
$ws.Range("I132").Value = 285715140
$ws.Range("J132").Value = 22003
$ws.Range("K132").Value = 857145420
$ws.Range("L132").Value = 66009
$ws.Range("M132").Value = -857142890
$ws.Range("N132").Value = -71069

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 9093203
$ws.Range("I132").Value = 15386035
$ws.Range("K132").Value = 46158105
$ws.Range("M132").Value = -46155575
$ws.Range("H136").Value = 5049.794
$ws.Range("I136").Value = 5538.857
$ws.Range("K136").Value = 16616.571
$ws.Range("M136").Value = -14066.571

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 85703880
$ws.Range("I132").Value = 120001410
$ws.Range("J132").Value = 51406356
$ws.Range("K132").Value = 360004230
$ws.Range("L132").Value = 154219068
$ws.Range("M132").Value = -360001700
$ws.Range("N132").Value = -154224128
$ws.Range("H136").Value = 25973528
$ws.Range("I136").Value = 17004640
$ws.Range("J136").Value = 41669084
$ws.Range("K136").Value = 51013920
$ws.Range("L136").Value = 125007252
$ws.Range("M136").Value = -51011370
$ws.Range("N136").Value = -125012352
